$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (contain only digits/dot) need to be
# forced to Text format first, otherwise Excel will silently convert the
# assigned string into a floating point number.
$numericLookingCells = @("D4","D5","D6","D8","D9","D11","D13","D14","D15","D17","D18","D19","D20","D21","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values coming from the refreshed crypto price feed.
$ws.Range("D2").Value = '29.935.88'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '1.891.92'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '0.8178'
$ws.Range("E5").Value = '  +6.28%  '
$ws.Range("D6").Value = '241.61'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").Value = '0.3222'
$ws.Range("E8").Value = '  +5.43%  '
$ws.Range("D9").Value = '26.40'
$ws.Range("E9").Value = '  +3.57%  '
$ws.Range("E10").Value = '  +2.54%  '
$ws.Range("D11").Value = '0.08035'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.903.02'
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.7451'
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").Value = '5.190'
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = '92.21'
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D16").Value = '29.946.88'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '14.00'
$ws.Range("E17").Value = '  +1.77%  '
$ws.Range("D18").Value = '5.888'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").Value = '244.24'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").Value = '0.000007743'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = '2.155.53'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").Value = '6.898'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '0.1580'
$ws.Range("E25").Value = '  +22.30%  '
$ws.Range("D26").Value = '166.08'
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").Value = '9.175'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("D29").Value = '2.067'
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("D30").Value = '1.368'
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").Value = '1.516'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").Value = '4.263'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").Value = '0.05596'
$ws.Range("E33").Value = '  +6.29%  '
$ws.Range("D34").Value = '4.070'
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").Value = '1.269'
$ws.Range("E35").Value = '  +1.96%  '
$ws.Range("D36").Value = '0.7307'
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").Value = '2.723'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").Value = '0.01911'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").Value = '2.790'
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").Value = '0.4408'
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").Value = '71.89'
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").Value = '5.945'
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("D43").Value = '0.8433'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").Value = '1.873'
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '100.65'
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.565'
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.679'
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '989.59'
$ws.Range("E49").Value = '  +8.56%  '
$ws.Range("D50").Value = '2.047.76'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '36.00'
$ws.Range("E51").Value = '  -0.69%  '

# Restore the default cell style for the cells we temporarily switched to
# Text format, so the final styling matches the original (unstyled) cells.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
